$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1443.0588
$ws.Range("I19").Value = 2347
$ws.Range("J19").Value = 1066.4166
$ws.Range("K19").Value = 2347
$ws.Range("L19").Value = 1066.4166
$ws.Range("M19").Value = -2172
$ws.Range("N19").Value = -1416.4166

$ws.Range("H39").Value = 102.166664
$ws.Range("I39").Value = 122.1
$ws.Range("J39").Value = 2.5
$ws.Range("K39").Value = 366.3
$ws.Range("L39").Value = 7.5
$ws.Range("M39").Value = -70.29999999999995
$ws.Range("N39").Value = -599.5

$ws.Range("H40").Value = 4150.636
$ws.Range("I40").Value = 2860.5
$ws.Range("J40").Value = 5698.8
$ws.Range("K40").Value = 2860.5
$ws.Range("L40").Value = 5698.8
$ws.Range("M40").Value = -2685.5
$ws.Range("N40").Value = -6048.8

$ws.Range("H53").Value = 91.818184
$ws.Range("I53").Value = 78.14286
$ws.Range("J53").Value = 115.75
$ws.Range("K53").Value = 78.14286
$ws.Range("L53").Value = 115.75
$ws.Range("M53").Value = 558.85714
$ws.Range("N53").Value = -1389.75

$ws.Range("H86").Value = 5249.125
$ws.Range("I86").Value = 3999.6667
$ws.Range("J86").Value = 8997.5
$ws.Range("K86").Value = 3999.6667
$ws.Range("L86").Value = 8997.5
$ws.Range("M86").Value = -2876.6667
$ws.Range("N86").Value = -11243.5

$ws.Range("H89").Value = 5249.125
$ws.Range("I89").Value = 3999.6667
$ws.Range("J89").Value = 8997.5
$ws.Range("K89").Value = 19998.3335
$ws.Range("L89").Value = 44987.5
$ws.Range("M89").Value = -14382.3335
$ws.Range("N89").Value = -56219.5

$ws.Range("H97").Value = 1756
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 1756
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 5268
$ws.Range("N97").Value = -6260

$ws.Range("H100").Value = 1997.25
$ws.Range("I100").Value = 1997.25
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 1997.25
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -1456.25

$ws.Range("H107").Value = 684.9167
$ws.Range("I107").Value = 561.3
$ws.Range("J107").Value = 1303
$ws.Range("K107").Value = 561.3
$ws.Range("L107").Value = 1303
$ws.Range("M107").Value = 1358.7
$ws.Range("N107").Value = -5143

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H21").Value = 5000
$ws.Range("I21").Value = 5000
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 5000
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = -4626
$ws.Range("N21").ClearContents()

$ws.Range("H32").Value = 3441.6316
$ws.Range("I32").Value = 3532.889
$ws.Range("J32").Value = 1799
$ws.Range("K32").Value = 3532.889
$ws.Range("L32").Value = 1799
$ws.Range("M32").Value = -3245.889
$ws.Range("N32").Value = -2373

$ws.Range("H43").Value = 31670.75
$ws.Range("I43").Value = 28671
$ws.Range("J43").Value = 34670.5
$ws.Range("K43").Value = 28671
$ws.Range("L43").Value = 34670.5
$ws.Range("M43").Value = -28358
$ws.Range("N43").Value = -35296.5

$ws.Range("H45").Value = 3540.182
$ws.Range("I45").Value = 1314
$ws.Range("J45").Value = 4375
$ws.Range("K45").Value = 1314
$ws.Range("L45").Value = 4375
$ws.Range("M45").Value = -937
$ws.Range("N45").Value = -5129

$ws.Range("H54").Value = 10000
$ws.Range("I54").Value = 10000
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 10000
$ws.Range("L54").Value = 0
$ws.Range("M54").Value = -9231

$ws.Range("H102").Value = 1599.5
$ws.Range("I102").Value = 1599
$ws.Range("J102").Value = 1600
$ws.Range("K102").Value = 1599
$ws.Range("L102").Value = 1600
$ws.Range("M102").Value = 23
$ws.Range("N102").Value = -4844

$ws.Range("H122").Value = 917.3333
$ws.Range("I122").Value = 917.3333
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 2751.9999
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -301.9998999999998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 930.6667
$ws.Range("I20").Value = 839.4
$ws.Range("J20").Value = 1044.75
$ws.Range("K20").Value = 839.4
$ws.Range("L20").Value = 1044.75
$ws.Range("M20").Value = -592.4
$ws.Range("N20").Value = -1538.75

$ws.Range("H86").Value = 7642.2856
$ws.Range("I86").Value = 2999.5
$ws.Range("J86").Value = 9499.4
$ws.Range("K86").Value = 2999.5
$ws.Range("L86").Value = 9499.4
$ws.Range("M86").Value = -1876.5
$ws.Range("N86").Value = -11745.4

$ws.Range("H89").Value = 7642.2856
$ws.Range("I89").Value = 2999.5
$ws.Range("J89").Value = 9499.4
$ws.Range("K89").Value = 14997.5
$ws.Range("L89").Value = 47497
$ws.Range("M89").Value = -9381.5
$ws.Range("N89").Value = -58729

$ws.Range("H99").Value = 2721.3
$ws.Range("I99").Value = 1961
$ws.Range("J99").Value = 3481.6
$ws.Range("K99").Value = 1961
$ws.Range("L99").Value = 3481.6
$ws.Range("M99").Value = -463
$ws.Range("N99").Value = -6477.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 10331.214
$ws.Range("I105").Value = 11063.7
$ws.Range("J105").Value = 8500
$ws.Range("K105").Value = 11063.7
$ws.Range("L105").Value = 8500
$ws.Range("M105").Value = -9316.700000000001
$ws.Range("N105").Value = -11994

$ws.Range("H107").Value = 1533.0588
$ws.Range("I107").Value = 1119.7693
$ws.Range("J107").Value = 2876.25
$ws.Range("K107").Value = 1119.7693
$ws.Range("L107").Value = 2876.25
$ws.Range("M107").Value = 800.2307000000001
$ws.Range("N107").Value = -6716.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H41").Value = 275
$ws.Range("I41").Value = 200
$ws.Range("J41").Value = 350
$ws.Range("K41").Value = 600
$ws.Range("L41").Value = 1050
$ws.Range("M41").Value = -262
$ws.Range("N41").Value = -1726

$ws.Range("H47").Value = 97
$ws.Range("I47").Value = 97
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 291
$ws.Range("L47").Value = 0
$ws.Range("M47").Value = 140

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 2500
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 2500
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 2500
$ws.Range("N7").Value = -2724

$ws.Range("H8").Value = 2500
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 2500
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 2500
$ws.Range("N8").Value = -2778

$ws.Range("H102").Value = 1600.8572
$ws.Range("I102").Value = 1651
$ws.Range("J102").Value = 1300
$ws.Range("K102").Value = 1651
$ws.Range("L102").Value = 1300
$ws.Range("M102").Value = -29
$ws.Range("N102").Value = -4544

$ws.Range("H126").Value = 4000
$ws.Range("I126").Value = 4000
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 12000
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -9530
$ws.Range("N126").ClearContents()

$ws.Range("H132").Value = 4049.6
$ws.Range("I132").Value = 3812.125
$ws.Range("J132").Value = 4999.5
$ws.Range("K132").Value = 11436.375
$ws.Range("L132").Value = 14998.5
$ws.Range("M132").Value = -8906.375
$ws.Range("N132").Value = -20058.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6000
$ws.Range("I7").Value = 6500
$ws.Range("J7").Value = 5000
$ws.Range("K7").Value = 6500
$ws.Range("L7").Value = 5000
$ws.Range("M7").Value = -6388
$ws.Range("N7").Value = -5224

$ws.Range("H22").Value = 9800
$ws.Range("I22").Value = 500
$ws.Range("J22").Value = 16000
$ws.Range("K22").Value = 500
$ws.Range("L22").Value = 16000
$ws.Range("M22").Value = -205
$ws.Range("N22").Value = -16590

$ws.Range("H27").Value = 9800
$ws.Range("I27").Value = 500
$ws.Range("J27").Value = 16000
$ws.Range("K27").Value = 500
$ws.Range("L27").Value = 16000
$ws.Range("M27").Value = -393
$ws.Range("N27").Value = -16214

$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("M40").ClearContents()

$ws.Range("H46").Value = 3543.111
$ws.Range("I46").Value = 3584.7144
$ws.Range("J46").Value = 3397.5
$ws.Range("K46").Value = 3584.7144
$ws.Range("L46").Value = 3397.5
$ws.Range("M46").Value = -3396.7144
$ws.Range("N46").Value = -3773.5

$ws.Range("H61").Value = 3499.5
$ws.Range("I61").Value = 3499.5
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 3499.5
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -3297.5

$ws.Range("H113").Value = 3499.5
$ws.Range("I113").Value = 3499.5
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 3499.5
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -1329.5

$ws.Range("H122").Value = 3183
$ws.Range("I122").Value = 2274.5
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 6823.5
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -4373.5
$ws.Range("N122").Value = -19900

$ws.Range("H126").Value = 6000
$ws.Range("I126").Value = 6500
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 19500
$ws.Range("L126").Value = 15000
$ws.Range("M126").Value = -17030
$ws.Range("N126").Value = -19940

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2667.6
$ws.Range("I132").Value = 2583.25
$ws.Range("J132").Value = 3005
$ws.Range("K132").Value = 7749.75
$ws.Range("L132").Value = 9015
$ws.Range("M132").Value = -5219.75
$ws.Range("N132").Value = -14075
